$wb = $excel.ActiveWorkbook
Write-Output $wb.CustomXMLParts.Count
for ($i=1; $i -le $wb.CustomXMLParts.Count; $i++) {
    $part = $wb.CustomXMLParts.Item($i)
    Write-Output $part.XML
}
